# AI architecture dead end
# Applies the sharedStrings / formatting rework described by the diff:
#  - tidies up a couple of note texts (typos / rewordings)
#  - adds two new notes ("Generic action scoring job?" on H2,
#    "CurrentAction: ... BlobRef" replacing the old C13 note)
#  - gives the "Response/Curve/Calculation" block (H6:H8 merge + H2/H9)
#    a new green highlight fill to set it apart from the grey blocks
#  - tightens up row 11 (drops the old wrapped note + its extra row height)
#  - makes C13 a wrapped, taller note cell
#  - moves the active selection to H5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Text content changes
# ---------------------------------------------------------------------

# B11: drop the second line of the note
$ws.Range("B11").Value = "component per action"

# C13: replace the old note with a new multi-line one
$ws.Range("C13").Value = "CurrentAction:" + [char]10 + "    BlobRef (maybe add int Id instead)"

# H2: brand new note cell
$ws.Range("H2").Value = "Generic action scoring job?"

# ---------------------------------------------------------------------
# 2. New "green" highlight fill for the Response/Curve/Calculation block
#    (H2, H6:H8 merged cell, H9) - mirrors the existing grey fill used
#    elsewhere (fillId 2 = theme accent3 @ 60% tint) but with accent6.
# ---------------------------------------------------------------------

$greenFill = 11854021   # RGB(197,224,180) == theme "accent6" tinted 60%, same recipe as the grey fill

foreach ($addr in @("H2", "H6", "H7", "H8", "H9")) {
    $r = $ws.Range($addr)
    $r.Interior.Color = $greenFill
    $r.Borders.LineStyle = 1
    $r.Borders.Weight = 2
}

# H6 keeps its centered + wrapped text, H7/H8 are centered (no wrap),
# H2/H9 keep default (left) alignment - only the fill/border are new there.
$ws.Range("H6").HorizontalAlignment = -4108   # xlCenter
$ws.Range("H6").WrapText = $true
$ws.Range("H7").HorizontalAlignment = -4108   # xlCenter
$ws.Range("H8").HorizontalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------------
# 3. Row 11: the note no longer wraps to two lines, so let the row
#    shrink back to the default height.
# ---------------------------------------------------------------------

$ws.Rows.Item(11).AutoFit()

# ---------------------------------------------------------------------
# 4. C13: wrap the new (longer) note and give its row more height.
# ---------------------------------------------------------------------

$ws.Range("C13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 45

# ---------------------------------------------------------------------
# 5. Selection moves from C16 to H5
# ---------------------------------------------------------------------

$ws.Range("H5").Select()
